# Revisão final dos diagramas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update textual tweaks: "Passo N" -> "(passo N)" style references
$ws.Range("B15").Value = " Alternativa 1 [Componente incompatível com outro já selecionado]  (passo 4)"
$ws.Range("B19").Value = " Alternativa 2 [Necessita de componentes extra] (passo 3)"
$ws.Range("B22").Value = " Exceção 3 [Utilizador rejeita alterações] (passos 3.2, 4.2 e 6)"

# 2. Remove the stray "s" value that lived in G9 (outside the real table)
$ws.Range("G9").ClearContents()

# 3. Row 9 has slightly shrunk now that G9 is empty
$ws.Rows.Item(9).RowHeight = 17.35

# 4. Move the active selection to where G9 used to be
$ws.Range("G9").Select()

# 5. Drop the stray trailing placeholder row left over at the sheet bottom
$ws.Rows.Item(1048576).Delete()
